# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" column (column E, rows 16-20) previously listed the
# contribution periods in descending order (1806, 1805, 1804, 1803, 1802).
# The database refresh re-sorts them in ascending order
# (1802, 1803, 1804, 1805, 1806) while everything else on the sheet
# (worker data, values, formatting) stays the same.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E16").Value = "1802"
$ws.Range("E17").Value = "1803"
$ws.Range("E18").Value = "1804"
$ws.Range("E19").Value = "1805"
$ws.Range("E20").Value = "1806"
